# Apply cryptos list price/volume update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.235.99"
$ws.Range("E2").Value = "'  -0.95%  "

$ws.Range("D3").Value = "'3.494.92"
$ws.Range("E3").Value = "'  -3.81%  "

$ws.Range("E4").Value = "'  +0.19%  "

$ws.Range("D5").Value = "'600.61"
$ws.Range("E5").Value = "'  -1.92%  "

$ws.Range("D6").Value = "'141.02"
$ws.Range("E6").Value = "'  -6.94%  "

$ws.Range("D7").Value = "'3.491.86"
$ws.Range("E7").Value = "'  -3.88%  "

$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "'  -0.20%  "

$ws.Range("E9").Value = "'  +2.93%  "

$ws.Range("D10").Value = "'7.58"
$ws.Range("E10").Value = "'  -5.31%  "

$ws.Range("D11").Value = "'0.127"
$ws.Range("E11").Value = "'  -7.95%  "

$ws.Range("D12").Value = "'0.398"
$ws.Range("E12").Value = "'  -4.68%  "

$ws.Range("D13").Value = "'4.093.50"
$ws.Range("E13").Value = "'  -3.56%  "

$ws.Range("D14").Value = "'0.0000190"
$ws.Range("E14").Value = "'  -9.53%  "

$ws.Range("D15").Value = "'28.09"
$ws.Range("E15").Value = "'  -7.38%  "

$ws.Range("D16").Value = "'3.506.74"
$ws.Range("E16").Value = "'  -3.50%  "

$ws.Range("E17").Value = "'  -0.92%  "

$ws.Range("D18").Value = "'66.274.23"
$ws.Range("E18").Value = "'  -1.02%  "

$ws.Range("D19").Value = "'10.52"
$ws.Range("E19").Value = "'  -10.66%  "

$ws.Range("D20").Value = "'6.01"
$ws.Range("E20").Value = "'  -6.31%  "

$ws.Range("D21").Value = "'14.44"
$ws.Range("E21").Value = "'  -4.86%  "

$ws.Range("D22").Value = "'415.82"
$ws.Range("E22").Value = "'  -3.14%  "

$ws.Range("D23").Value = "'0.577"
$ws.Range("E23").Value = "'  -6.82%  "

$ws.Range("D24").Value = "'76.24"
$ws.Range("E24").Value = "'  -3.65%  "

$ws.Range("D25").Value = "'3.648.16"
$ws.Range("E25").Value = "'  -3.17%  "

$ws.Range("E26").Value = "'  -0.09%  "

$ws.Range("D27").Value = "'0.0000109"
$ws.Range("E27").Value = "'  -12.43%  "

$ws.Range("D28").Value = "'2.41"
$ws.Range("E28").Value = "'  -4.75%  "

$ws.Range("D29").Value = "'8.76"
$ws.Range("E29").Value = "'  -9.12%  "

$ws.Range("D30").Value = "'7.56"
$ws.Range("E30").Value = "'  -10.65%  "

$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "'  +0.00%  "

$ws.Range("D32").Value = "'3.510.97"
$ws.Range("E32").Value = "'  -3.24%  "

$ws.Range("D33").Value = "'0.151"
$ws.Range("E33").Value = "'  -6.37%  "

$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "'  +0.03%  "

$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = "'23.78"
$ws.Range("E35").Value = "'  -6.85%  "

$ws.Range("D36").Value = "'1.30"
$ws.Range("E36").Value = "'  -13.08%  "

$ws.Range("D37").Value = "'7.30"
$ws.Range("E37").Value = "'  -7.78%  "

$ws.Range("D38").Value = "'1.60"
$ws.Range("E38").Value = "'  -6.96%  "

$ws.Range("D39").Value = "'173.30"
$ws.Range("E39").Value = "'  -2.54%  "

$ws.Range("D40").Value = "'5.06"
$ws.Range("E40").Value = "'  -11.48%  "

$ws.Range("D41").Value = "'0.0790"
$ws.Range("E41").Value = "'  -8.78%  "

$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").Value = "'0.842"
$ws.Range("E42").Value = "'  -6.68%  "

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = "'4.84"
$ws.Range("E43").Value = "'  -7.95%  "

$ws.Range("D44").Value = "'45.05"
$ws.Range("E44").Value = "'  -2.66%  "

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "'  +0.03%  "

$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").Value = "'1.70"
$ws.Range("E46").Value = "'  -11.13%  "

$ws.Range("B47").Value = 'Cosmos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D47").Value = "'6.92"
$ws.Range("E47").Value = "'  -4.39%  "

$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").Value = "'2.27"
$ws.Range("E48").Value = "'  -13.46%  "

$ws.Range("E49").Value = "'  -6.47%  "

$ws.Range("D50").Value = "'22.40"
$ws.Range("E50").Value = "'  -7.50%  "

$ws.Range("D51").Value = "'0.872"
$ws.Range("E51").Value = "'  -11.27%  "

Write-Output "Updated cryptos list"